$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03739799725335
$ws.Cells.Item(2, 4).Value = 1.040042358767425
$ws.Cells.Item(2, 5).Value = 1.045425107642324
$ws.Cells.Item(2, 6).Value = 1.054492432429486
$ws.Cells.Item(2, 9).Value = 1.035492039141551
$ws.Cells.Item(2, 10).Value = 1.042501087052323
$ws.Cells.Item(2, 11).Value = 1.042825830249702
$ws.Cells.Item(2, 12).Value = 1.048193384254906
$ws.Cells.Item(2, 13).Value = 1.05723549321273
$ws.Cells.Item(2, 14).Value = 1.043981558172484
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038416890594237
$ws.Cells.Item(3, 4).Value = 1.040975615978251
$ws.Cells.Item(3, 5).Value = 1.046335954595482
$ws.Cells.Item(3, 6).Value = 1.055499945939038
$ws.Cells.Item(3, 9).Value = 1.035641983404008
$ws.Cells.Item(3, 10).Value = 1.043163994072834
$ws.Cells.Item(3, 11).Value = 1.043569106176456
$ws.Cells.Item(3, 12).Value = 1.048915410179594
$ws.Cells.Item(3, 13).Value = 1.058055764217448
$ws.Cells.Item(3, 14).Value = 1.044645406597001
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.039076488622205
$ws.Cells.Item(4, 4).Value = 1.041580070794917
$ws.Cells.Item(4, 5).Value = 1.046925965896413
$ws.Cells.Item(4, 6).Value = 1.056152582956061
$ws.Cells.Item(4, 9).Value = 1.035737620874194
$ws.Cells.Item(4, 10).Value = 1.043592657789618
$ws.Cells.Item(4, 11).Value = 1.044050005489704
$ws.Cells.Item(4, 12).Value = 1.049382610373087
$ws.Cells.Item(4, 13).Value = 1.058586629840716
$ws.Cells.Item(4, 14).Value = 1.04507467906543
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.039353856464087
$ws.Cells.Item(5, 4).Value = 1.041834320399949
$ws.Cells.Item(5, 5).Value = 1.047174157149689
$ws.Cells.Item(5, 6).Value = 1.056427120269508
$ws.Cells.Item(5, 9).Value = 1.035777494544909
$ws.Cells.Item(5, 10).Value = 1.043772800000639
$ws.Cells.Item(5, 11).Value = 1.044252162907254
$ws.Cells.Item(5, 12).Value = 1.049579020896009
$ws.Cells.Item(5, 13).Value = 1.058809827808281
$ws.Cells.Item(5, 14).Value = 1.045255077099053
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.039400432004439
$ws.Cells.Item(6, 4).Value = 1.041877018024588
$ws.Cells.Item(6, 5).Value = 1.047215838352628
$ws.Cells.Item(6, 6).Value = 1.056473226160184
$ws.Cells.Item(6, 9).Value = 1.035784170010782
$ws.Cells.Item(6, 10).Value = 1.043803042663639
$ws.Cells.Item(6, 11).Value = 1.044286105264727
$ws.Cells.Item(6, 12).Value = 1.0496119990407
$ws.Cells.Item(6, 13).Value = 1.058847304999886
$ws.Cells.Item(6, 14).Value = 1.045285362710104
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.039080194539202
$ws.Cells.Item(7, 4).Value = 1.041583467554121
$ws.Cells.Item(7, 5).Value = 1.04692928164984
$ws.Cells.Item(7, 6).Value = 1.056156250675261
$ws.Cells.Item(7, 9).Value = 1.035738154974418
$ws.Cells.Item(7, 10).Value = 1.043595065126707
$ws.Cells.Item(7, 11).Value = 1.044052706777613
$ws.Cells.Item(7, 12).Value = 1.04938523482272
$ws.Cells.Item(7, 13).Value = 1.058589612137189
$ws.Cells.Item(7, 14).Value = 1.045077089821214
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037742273445193
$ws.Cells.Item(8, 4).Value = 1.04035763784025
$ws.Cells.Item(8, 5).Value = 1.045732801148134
$ws.Cells.Item(8, 6).Value = 1.054832779126246
$ws.Cells.Item(8, 9).Value = 1.035543000353892
$ws.Cells.Item(8, 10).Value = 1.042725177478257
$ws.Cells.Item(8, 11).Value = 1.043077033740188
$ws.Cells.Item(8, 12).Value = 1.048437395632903
$ws.Cells.Item(8, 13).Value = 1.057512687405536
$ws.Cells.Item(8, 14).Value = 1.044205966832526
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.035387045983139
$ws.Cells.Item(9, 4).Value = 1.038202010616398
$ws.Cells.Item(9, 5).Value = 1.043629332710534
$ws.Cells.Item(9, 6).Value = 1.052506117483824
$ws.Cells.Item(9, 9).Value = 1.035188511350434
$ws.Cells.Item(9, 10).Value = 1.041190196577821
$ws.Cells.Item(9, 11).Value = 1.041357421845094
$ws.Cells.Item(9, 12).Value = 1.046767221845218
$ws.Cells.Item(9, 13).Value = 1.055615773932045
$ws.Cells.Item(9, 14).Value = 1.042668806083146
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033818502765769
$ws.Cells.Item(10, 4).Value = 1.036767960942178
$ws.Cells.Item(10, 5).Value = 1.04223035531598
$ws.Cells.Item(10, 6).Value = 1.050958734188874
$ws.Cells.Item(10, 9).Value = 1.034945080151191
$ws.Cells.Item(10, 10).Value = 1.040165481606809
$ws.Cells.Item(10, 11).Value = 1.040210817924699
$ws.Cells.Item(10, 12).Value = 1.045653840264502
$ws.Cells.Item(10, 13).Value = 1.054351724743134
$ws.Cells.Item(10, 14).Value = 1.041642635899339
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033139691484185
$ws.Cells.Item(11, 4).Value = 1.036147731309878
$ws.Cells.Item(11, 5).Value = 1.041625384090438
$ws.Cells.Item(11, 6).Value = 1.05028959363901
$ws.Cells.Item(11, 9).Value = 1.034837991100663
$ws.Cells.Item(11, 10).Value = 1.039721445159279
$ws.Cells.Item(11, 11).Value = 1.03971428728705
$ws.Cells.Item(11, 12).Value = 1.045171757925326
$ws.Cells.Item(11, 13).Value = 1.053804519279505
$ws.Cells.Item(11, 14).Value = 1.041197968869122
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032887607791049
$ws.Cells.Item(12, 4).Value = 1.035917459812632
$ws.Cells.Item(12, 5).Value = 1.041400791099917
$ws.Cells.Item(12, 6).Value = 1.050041179001774
$ws.Cells.Item(12, 9).Value = 1.034797961129166
$ws.Cells.Item(12, 10).Value = 1.039556461475485
$ws.Cells.Item(12, 11).Value = 1.039529847737588
$ws.Cells.Item(12, 12).Value = 1.044992694568289
$ws.Cells.Item(12, 13).Value = 1.053601283928883
$ws.Cells.Item(12, 14).Value = 1.041032750889573
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032941678014334
$ws.Cells.Item(13, 4).Value = 1.035966848867973
$ws.Cells.Item(13, 5).Value = 1.041448961604177
$ws.Cells.Item(13, 6).Value = 1.050094458710843
$ws.Cells.Item(13, 9).Value = 1.034806559107066
$ws.Cells.Item(13, 10).Value = 1.039591853234629
$ws.Cells.Item(13, 11).Value = 1.039569410907478
$ws.Cells.Item(13, 12).Value = 1.045031104098732
$ws.Cells.Item(13, 13).Value = 1.053644877634961
$ws.Cells.Item(13, 14).Value = 1.041068192909074
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033118853002633
$ws.Cells.Item(14, 4).Value = 1.036128694760048
$ws.Cells.Item(14, 5).Value = 1.041606816705981
$ws.Cells.Item(14, 6).Value = 1.050269056870583
$ws.Cells.Item(14, 9).Value = 1.034834687356265
$ws.Cells.Item(14, 10).Value = 1.039707808548538
$ws.Cells.Item(14, 11).Value = 1.03969904157461
$ws.Cells.Item(14, 12).Value = 1.045156956417743
$ws.Cells.Item(14, 13).Value = 1.053787719347277
$ws.Cells.Item(14, 14).Value = 1.04118431289283
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033228023925227
$ws.Cells.Item(15, 4).Value = 1.03622842787469
$ws.Cells.Item(15, 5).Value = 1.041704092403102
$ws.Cells.Item(15, 6).Value = 1.050376650332286
$ws.Cells.Item(15, 9).Value = 1.034851984672712
$ws.Cells.Item(15, 10).Value = 1.039779245996591
$ws.Cells.Item(15, 11).Value = 1.039778910534366
$ws.Cells.Item(15, 12).Value = 1.045234498665328
$ws.Cells.Item(15, 13).Value = 1.053875731659482
$ws.Cells.Item(15, 14).Value = 1.041255851790253
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.033863561222562
$ws.Cells.Item(16, 4).Value = 1.036809138842445
$ws.Cells.Item(16, 5).Value = 1.042270522085239
$ws.Cells.Item(16, 6).Value = 1.051003161587671
$ws.Cells.Item(16, 9).Value = 1.034952151894851
$ws.Cells.Item(16, 10).Value = 1.040194943974805
$ws.Cells.Item(16, 11).Value = 1.040243770159594
$ws.Cells.Item(16, 12).Value = 1.045685834937961
$ws.Cells.Item(16, 13).Value = 1.054388043893295
$ws.Cells.Item(16, 14).Value = 1.041672140107277
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034262318401207
$ws.Cells.Item(17, 4).Value = 1.037173597664266
$ws.Cells.Item(17, 5).Value = 1.04262604191751
$ws.Cells.Item(17, 6).Value = 1.051396393450924
$ws.Cells.Item(17, 9).Value = 1.035014534132605
$ws.Cells.Item(17, 10).Value = 1.04045561285635
$ws.Cells.Item(17, 11).Value = 1.040535353067043
$ws.Cells.Item(17, 12).Value = 1.045968951829335
$ws.Cells.Item(17, 13).Value = 1.054709440749118
$ws.Cells.Item(17, 14).Value = 1.041933179168534
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034494943201634
$ws.Cells.Item(18, 4).Value = 1.037386250071668
$ws.Cells.Item(18, 5).Value = 1.042833487295002
$ws.Cells.Item(18, 6).Value = 1.051625844440122
$ws.Cells.Item(18, 9).Value = 1.035050758351142
$ws.Cells.Item(18, 10).Value = 1.040607624912848
$ws.Cells.Item(18, 11).Value = 1.04070542420032
$ws.Cells.Item(18, 12).Value = 1.046134090937545
$ws.Cells.Item(18, 13).Value = 1.054896919191983
$ws.Cells.Item(18, 14).Value = 1.042085407099587
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034574268429158
$ws.Cells.Item(19, 4).Value = 1.037458770854546
$ws.Cells.Item(19, 5).Value = 1.042904233789223
$ws.Cells.Item(19, 6).Value = 1.051704095776741
$ws.Cells.Item(19, 9).Value = 1.035063082337064
$ws.Cells.Item(19, 10).Value = 1.040659451704548
$ws.Cells.Item(19, 11).Value = 1.040763413324229
$ws.Cells.Item(19, 12).Value = 1.046190399395589
$ws.Cells.Item(19, 13).Value = 1.054960846696996
$ws.Cells.Item(19, 14).Value = 1.042137307491276
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.034219531754375
$ws.Cells.Item(20, 4).Value = 1.037134487462726
$ws.Cells.Item(20, 5).Value = 1.042587890067598
$ws.Cells.Item(20, 6).Value = 1.051354194544733
$ws.Cells.Item(20, 9).Value = 1.035007857892391
$ws.Cells.Item(20, 10).Value = 1.04042764882656
$ws.Cells.Item(20, 11).Value = 1.04050406941523
$ws.Cells.Item(20, 12).Value = 1.045938575874122
$ws.Cells.Item(20, 13).Value = 1.05467495654888
$ws.Cells.Item(20, 14).Value = 1.041905175426614
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033066677808245
$ws.Cells.Item(21, 4).Value = 1.03608103215278
$ws.Cells.Item(21, 5).Value = 1.041560328986877
$ws.Cells.Item(21, 6).Value = 1.050217638360392
$ws.Cells.Item(21, 9).Value = 1.034826411250829
$ws.Cells.Item(21, 10).Value = 1.039673663941009
$ws.Cells.Item(21, 11).Value = 1.039660868734
$ws.Cells.Item(21, 12).Value = 1.045119895945466
$ws.Cells.Item(21, 13).Value = 1.053745655434801
$ws.Cells.Item(21, 14).Value = 1.041150119796042
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032342162903357
$ws.Cells.Item(22, 4).Value = 1.035419316327604
$ws.Cells.Item(22, 5).Value = 1.040914956723012
$ws.Cells.Item(22, 6).Value = 1.0495038157813
$ws.Cells.Item(22, 9).Value = 1.034710868807093
$ws.Cells.Item(22, 10).Value = 1.039199321769424
$ws.Cells.Item(22, 11).Value = 1.039130680791323
$ws.Cells.Item(22, 12).Value = 1.044605179861464
$ws.Cells.Item(22, 13).Value = 1.053161489480898
$ws.Cells.Item(22, 14).Value = 1.040675104004165
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032726210599705
$ws.Cells.Item(23, 4).Value = 1.035770044129192
$ws.Cells.Item(23, 5).Value = 1.041257014519139
$ws.Cells.Item(23, 6).Value = 1.049882152886022
$ws.Cells.Item(23, 9).Value = 1.034772258298006
$ws.Cells.Item(23, 10).Value = 1.039450806081833
$ws.Cells.Item(23, 11).Value = 1.039411746498472
$ws.Cells.Item(23, 12).Value = 1.044878038454755
$ws.Cells.Item(23, 13).Value = 1.053471155068371
$ws.Cells.Item(23, 14).Value = 1.040926945453143
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034238865085522
$ws.Cells.Item(24, 4).Value = 1.037152159463837
$ws.Cells.Item(24, 5).Value = 1.042605129009457
$ws.Cells.Item(24, 6).Value = 1.051373262149844
$ws.Cells.Item(24, 9).Value = 1.035010875099412
$ws.Cells.Item(24, 10).Value = 1.040440284665835
$ws.Cells.Item(24, 11).Value = 1.040518205164
$ws.Cells.Item(24, 12).Value = 1.045952301455155
$ws.Cells.Item(24, 13).Value = 1.054690538434152
$ws.Cells.Item(24, 14).Value = 1.041917829210229
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035995646820954
$ws.Cells.Item(25, 4).Value = 1.038758760551486
$ws.Cells.Item(25, 5).Value = 1.044172545716439
$ws.Cells.Item(25, 6).Value = 1.053106962532558
$ws.Cells.Item(25, 9).Value = 1.035281409261988
$ws.Cells.Item(25, 10).Value = 1.041587274353735
$ws.Cells.Item(25, 11).Value = 1.041802020231986
$ws.Cells.Item(25, 12).Value = 1.047198992871349
$ws.Cells.Item(25, 13).Value = 1.056106076173236
$ws.Cells.Item(25, 14).Value = 1.043066447755047
